$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns retain their original text formatting so that
# numeric-looking strings (e.g. "241.70", "0.7184") are not coerced
# into Double values by Excel, matching the original inlineStr cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.475.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7184"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07892"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3098"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.48"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08272"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7285"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.498.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "245.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007856"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.130.07"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.084"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1619"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +14.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.74%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.354"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.97%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.400"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.108"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05203"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.198"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7271"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.204.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.697"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9101"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.177"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.34"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.025.74"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5291"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.800"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.926"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.72%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4326"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.68%  "
